$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 9.75
$ws.Range("H5").Value = 5.5
$ws.Range("I5").Value = 1.21
$ws.Range("T5").Value = 24
$ws.Range("U5").Value = 60
$ws.Range("V5").Value = 26
$ws.Range("W5").Value = 200
$ws.Range("X5").Value = 90
$ws.Range("Y5").Value = 70
$ws.Range("Z5").Value = 17
$ws.Range("AA5").Value = 10.25
$ws.Range("AB5").Value = 19.5
$ws.Range("AC5").Value = 75
$ws.Range("AD5").Value = 450
$ws.Range("AF5").Value = 5.9
$ws.Range("AG5").Value = 7.9
$ws.Range("AH5").Value = 6.4
$ws.Range("AJ5").Value = 21

# Row 7
$ws.Range("G7").Value = 1.55
$ws.Range("H7").Value = 3.95
$ws.Range("I7").Value = 4.7
$ws.Range("U7").Value = 8
$ws.Range("V7").Value = 7
$ws.Range("W7").Value = 10.75
$ws.Range("Y7").Value = 15
$ws.Range("Z7").Value = 16.5
$ws.Range("AA7").Value = 7.4
$ws.Range("AB7").Value = 11
$ws.Range("AE7").Value = 16
$ws.Range("AF7").Value = 27
$ws.Range("AG7").Value = 13
$ws.Range("AJ7").Value = 27

# Row 10
$ws.Range("L10").Value = 1.3
$ws.Range("M10").Value = 3.4
$ws.Range("Z10").Value = 9.5

# Row 13
$ws.Range("R13").Value = 1.7
